$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 00:34:21"
$wsZhCn.Range("H3").Value = "2016-03-18 00:34:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 00:34:24"
$wsDeDe.Range("H3").Value = "2016-03-18 00:34:46"
